# Apply the "1st iteration" update to the Metadata worksheet:
#  - set the "Experimental" value (B7) to "true"
#  - refresh the "Date" value (B8) to the new generation timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B7").Value = "true"
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
